$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Free up the old H17/I17 shared strings and replace with the new Siege Tank values
# (this keeps "Siege mode"/"Siege Tank" referenced while releasing
#  "damage/speed/freq/range" and "330/-4/-4/4")
$ws.Range("H17").Value = "l/h/speed/freq/range"
$ws.Range("I17").Value = "310/330/-4/-4/4"

# Insert a new row above the Siege Tank row, pushing it down to row 18
$ws.Rows.Item(17).Insert()

# Give the new blank row 17 the same formatting as the row above (row 16)
$ws.Range("A16:J16").Copy()
$ws.Paste($ws.Range("A17:J17"))

# Fill in the new Storm Chariot / EMP skill row
$ws.Range("A17").Value = "Storm Chariot"
$ws.Range("B17").Value = "EMP"
$ws.Range("C17").Value = "continuous"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = "n"
$ws.Range("F17").Value = 75
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = "shield/speed"
$ws.Range("I17").Value = "-50/2"
$ws.Range("J17").Value = 6

# Add new rows 19-21 using row 18 (Siege Tank, switch type) as a formatting template
$ws.Range("A18:J18").Copy()
$ws.Paste($ws.Range("A19:J19"))
$ws.Paste($ws.Range("A20:J20"))
$ws.Paste($ws.Range("A21:J21"))

# Row 19: Frigate / Piercing missile
$ws.Range("A19").Value = "Frigate"
$ws.Range("B19").Value = "Piercing missile"
$ws.Range("C19").Value = "damage"
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = "g_m"
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = "n"
$ws.Range("I19").Value = "500"
$ws.Range("J19").Value = 8

# Row 20: Battleship / Space jump
$ws.Range("A20").Value = "Battleship"
$ws.Range("B20").Value = "Space jump"
$ws.Range("C20").Value = "continuous"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = "self"
$ws.Range("F20").Value = 150
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = "speed"
$ws.Range("I20").Value = "10"
$ws.Range("J20").Value = 0

# Row 21: Battleship / Volley
$ws.Range("A21").Value = "Battleship"
$ws.Range("B21").Value = "Volley"
$ws.Range("C21").Value = "continuous"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = "self"
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = "g/a/g_freq/a_freq"
$ws.Range("I21").Value = "550/180/-14/-34"
$ws.Range("J21").Value = 0

# Column width adjustments
$ws.Columns.Item(4).ColumnWidth = 4.8
$ws.Columns.Item(5).ColumnWidth = 6.8
$ws.Columns.Item(6).ColumnWidth = 7.5
$ws.Columns.Item(7).ColumnWidth = 5.5

# Selection
$ws.Range("L15").Select()
